$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "sat_1"
$ws.Range("B3").Value = "sat_2"
$ws.Range("B4").Value = "sat_3"
$ws.Range("B5").Value = "sat_4"
$ws.Range("B6").Value = "sat_5"
$ws.Range("B7").Value = "sat_6"
$ws.Range("B8").Value = "sat_7"

[void]$ws.Range("B2:B8").Select()

